$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings such as "42.990.86" (dotted thousands-style)
# or plain decimals such as "1.01". Force each target cell to Text format
# before writing so Excel does not auto-coerce the plain-decimal-looking
# ones into numeric cells (which would lose the exact printed text / trailing
# zeros and introduce floating point noise).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.990.86"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.306.29"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.03"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.82"
$ws.Range("E6").Value = "  +2.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.605"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.55"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0906"
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.36"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.108"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.993"
$ws.Range("E14").Value = "  +3.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.23"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.656.36"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.307.50"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.854.64"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.34"
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000105"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.55"
$ws.Range("E21").Value = "  +4.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.40"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.53"
$ws.Range("E23").Value = "  -1.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.98"
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("E25").Value = "  -1.81%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.78"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.33"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.88"
$ws.Range("E29").Value = "  +14.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.46"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.34"
$ws.Range("E31").Value = "  -2.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.94"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0868"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.131"
$ws.Range("E34").Value = "  -1.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.56"
$ws.Range("E35").Value = "  -3.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.114"
$ws.Range("E36").Value = "  -2.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.53"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0351"
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.77"
$ws.Range("E39").Value = "  +2.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.69"
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.59"
$ws.Range("E41").Value = "  +3.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.53"
$ws.Range("E42").Value = "  +4.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.07"
$ws.Range("E43").Value = "  +1.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.228"
$ws.Range("E44").Value = "  +1.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.19"
$ws.Range("E46").Value = "  +1.70%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.747.36"
$ws.Range("E47").Value = "  +8.93%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "81.25"
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "111.93"
$ws.Range("E49").Value = "  -2.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.23"
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.69"
$ws.Range("E51").Value = "  -2.73%  "
